$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column C (rows 2-31), after re-sorting the data by ratio
# the feasible/best values need to be recalculated, so column C is updated
# with the new set of numbers.
$newValues = @{
    2  = 23922
    3  = 24717
    4  = 24623
    5  = 23313
    6  = 24709
    7  = 25496
    8  = 25470
    9  = 23510
    10 = 24385
    11 = 24368
    12 = 41955
    13 = 41201
    14 = 41729
    15 = 43021
    16 = 40829
    17 = 41028
    18 = 40077
    19 = 43892
    20 = 43449
    21 = 44177
    22 = 60079
    23 = 62160
    24 = 59439
    25 = 60077
    26 = 60468
    27 = 59430
    28 = 61439
    29 = 61071
    30 = 58841
    31 = 59821
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 3).Value = $newValues[$row]
}

$excel.Calculate()

# Update the selected cell/range to match the saved view state
$ws.Range("G6").Select() | Out-Null
